# Auto-generated Excel COM-interop script to apply the commit diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 32000
$ws.Range("I9").Value = 45500
$ws.Range("J9").Value = 5000
$ws.Range("K9").Value = 45500
$ws.Range("L9").Value = 5000
$ws.Range("M9").Value = -45331
$ws.Range("N9").Value = -5338

$ws.Range("H11").Value = 537.5
$ws.Range("I11").Value = 537.5
$ws.Range("K11").Value = 537.5
$ws.Range("M11").Value = -397.5

$ws.Range("H40").Value = 2388.6
$ws.Range("I40").Value = 1475
$ws.Range("K40").Value = 1475
$ws.Range("M40").Value = -1300

$ws.Range("H41").Value = 592.4737
$ws.Range("I41").Value = 556.0833
$ws.Range("J41").Value = 654.8570999999999
$ws.Range("K41").Value = 556.0833
$ws.Range("L41").Value = 654.8570999999999
$ws.Range("M41").Value = -116.0833
$ws.Range("N41").Value = -1534.8571

$ws.Range("H96").Value = 2088.3333
$ws.Range("I96").Value = 1927.25
$ws.Range("J96").Value = 2217.2
$ws.Range("K96").Value = 5781.75
$ws.Range("L96").Value = 6651.599999999999
$ws.Range("M96").Value = -4408.75
$ws.Range("N96").Value = -9397.599999999999

$ws.Range("H129").Value = 2093.7778
$ws.Range("I129").Value = 1953.375
$ws.Range("K129").Value = 5860.125
$ws.Range("M129").Value = -860.125

$ws.Range("H135").Value = 1170.76
$ws.Range("I135").Value = 1031.7368
$ws.Range("K135").Value = 9285.6312
$ws.Range("M135").Value = -6750.6312

$ws.Range("H138").Value = 10873350
$ws.Range("I138").Value = 1298
$ws.Range("K138").Value = 3894
$ws.Range("M138").Value = 1246

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 1755.5769
$ws.Range("I88").Value = 918.8333
$ws.Range("J88").Value = 2006.6
$ws.Range("K88").Value = 918.8333
$ws.Range("L88").Value = 2006.6
$ws.Range("M88").Value = -512.8333
$ws.Range("N88").Value = -2818.6

$ws.Range("H91").Value = 1755.5769
$ws.Range("I91").Value = 918.8333
$ws.Range("J91").Value = 2006.6
$ws.Range("K91").Value = 918.8333
$ws.Range("L91").Value = 2006.6
$ws.Range("M91").Value = 485.1667
$ws.Range("N91").Value = -4814.6

$ws.Range("H105").Value = 40000
$ws.Range("J105").Value = 40000
$ws.Range("L105").Value = 40000
$ws.Range("N105").Value = -46988

$ws.Range("H132").Value = 5293.7
$ws.Range("I132").Value = 5104.5186
$ws.Range("K132").Value = 15313.5558
$ws.Range("M132").Value = -12783.5558

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H69").Value = 45000
$ws.Range("J69").Value = 45000
$ws.Range("L69").Value = 45000
$ws.Range("N69").Value = -46622

$ws.Range("H72").Value = 45000
$ws.Range("J72").Value = 45000
$ws.Range("L72").Value = 135000
$ws.Range("N72").Value = -143112

$ws.Range("H86").Value = 2586.5625
$ws.Range("I86").Value = 2470.9333
$ws.Range("J86").Value = 4321
$ws.Range("K86").Value = 2470.9333
$ws.Range("L86").Value = 4321
$ws.Range("M86").Value = -1347.9333
$ws.Range("N86").Value = -6567

$ws.Range("H89").Value = 2586.5625
$ws.Range("I89").Value = 2470.9333
$ws.Range("J89").Value = 4321
$ws.Range("K89").Value = 12354.6665
$ws.Range("L89").Value = 21605
$ws.Range("M89").Value = -6738.666500000001
$ws.Range("N89").Value = -32837

$ws.Range("H100").Value = 17657.166
$ws.Range("J100").Value = 17657.166
$ws.Range("L100").Value = 17657.166
$ws.Range("N100").Value = -19821.166

$ws.Range("H134").Value = 1375.1772
$ws.Range("I134").Value = 1367.1666
$ws.Range("K134").Value = 4101.4998
$ws.Range("M134").Value = -1566.4998

$ws.Range("H135").Value = 78780
$ws.Range("J135").Value = 78780
$ws.Range("L135").Value = 78780
$ws.Range("N135").Value = -88920

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 981.4666999999999
$ws.Range("J7").Value = 650.6667
$ws.Range("L7").Value = 650.6667
$ws.Range("N7").Value = -876.6667

$ws.Range("H43").Value = 29832.666
$ws.Range("J43").Value = 29832.666
$ws.Range("L43").Value = 29832.666
$ws.Range("N43").Value = -30200.666

$ws.Range("H62").Value = 7452.273
$ws.Range("I62").Value = 6197
$ws.Range("J62").Value = 8498.333000000001
$ws.Range("K62").Value = 6197
$ws.Range("L62").Value = 8498.333000000001
$ws.Range("M62").Value = -5573
$ws.Range("N62").Value = -9746.333000000001

$ws.Range("H65").Value = 7452.273
$ws.Range("I65").Value = 6197
$ws.Range("J65").Value = 8498.333000000001
$ws.Range("K65").Value = 30985
$ws.Range("L65").Value = 42491.665
$ws.Range("M65").Value = -27865
$ws.Range("N65").Value = -48731.665

$ws.Range("H101").Value = 29832.666
$ws.Range("J101").Value = 29832.666
$ws.Range("L101").Value = 29832.666
$ws.Range("N101").Value = -36322.666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 5003.8887
$ws.Range("I62").Value = 2758.75
$ws.Range("J62").Value = 6800
$ws.Range("K62").Value = 8276.25
$ws.Range("L62").Value = 20400
$ws.Range("M62").Value = -7590.25
$ws.Range("N62").Value = -21772

$ws.Range("H65").Value = 5003.8887
$ws.Range("I65").Value = 2758.75
$ws.Range("J65").Value = 6800
$ws.Range("K65").Value = 24828.75
$ws.Range("L65").Value = 61200
$ws.Range("M65").Value = -21396.75
$ws.Range("N65").Value = -68064

$ws.Range("H80").Value = 7249.75
$ws.Range("I80").Value = 4500
$ws.Range("K80").Value = 13500
$ws.Range("M80").Value = -12564

$ws.Range("H83").Value = 7249.75
$ws.Range("I83").Value = 4500
$ws.Range("K83").Value = 40500
$ws.Range("M83").Value = -35820

$ws.Range("H113").Value = 748.25
$ws.Range("I113").Value = 900
$ws.Range("J113").Value = 726.5714
$ws.Range("K113").Value = 2700
$ws.Range("L113").Value = 2179.7142
$ws.Range("M113").Value = -530
$ws.Range("N113").Value = -6519.7142

$ws.Range("H122").Value = 736.3333
$ws.Range("I122").Value = 604
$ws.Range("J122").Value = 802.5
$ws.Range("K122").Value = 5436
$ws.Range("L122").Value = 7222.5
$ws.Range("M122").Value = -2986
$ws.Range("N122").Value = -12122.5

$ws.Range("H129").Value = 1421.9
$ws.Range("I129").Value = 1027.375
$ws.Range("K129").Value = 3082.125
$ws.Range("M129").Value = 1917.875

$ws.Range("H131").Value = 37174.07
$ws.Range("I131").Value = 91811.55
$ws.Range("J131").Value = 1820.4117
$ws.Range("K131").Value = 275434.65
$ws.Range("L131").Value = 5461.2351
$ws.Range("M131").Value = -270394.65
$ws.Range("N131").Value = -15541.2351

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3207.9285
$ws.Range("I122").Value = 3665.125
$ws.Range("J122").Value = 2598.3333
$ws.Range("K122").Value = 10995.375
$ws.Range("L122").Value = 7794.999899999999
$ws.Range("M122").Value = -8545.375
$ws.Range("N122").Value = -12694.9999

$ws.Range("H132").Value = 9151.75
$ws.Range("I132").Value = 9737.706
$ws.Range("J132").Value = 5831.3335
$ws.Range("K132").Value = 29213.118
$ws.Range("L132").Value = 17494.0005
$ws.Range("M132").Value = -26683.118
$ws.Range("N132").Value = -22554.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6142.2666
$ws.Range("I7").Value = 5784.1113
$ws.Range("J7").Value = 6679.5
$ws.Range("K7").Value = 5784.1113
$ws.Range("L7").Value = 6679.5
$ws.Range("M7").Value = -5672.1113
$ws.Range("N7").Value = -6903.5

$ws.Range("H16").Value = 1155.5862
$ws.Range("I16").Value = 1048.3182
$ws.Range("J16").Value = 1492.7142
$ws.Range("K16").Value = 1048.3182
$ws.Range("L16").Value = 1492.7142
$ws.Range("M16").Value = -878.3181999999999
$ws.Range("N16").Value = -1832.7142

$ws.Range("H22").Value = 1772.1177
$ws.Range("I22").Value = 1759.3334
$ws.Range("J22").Value = 1802.8
$ws.Range("K22").Value = 1759.3334
$ws.Range("L22").Value = 1802.8
$ws.Range("M22").Value = -1464.3334
$ws.Range("N22").Value = -2392.8

$ws.Range("H27").Value = 1772.1177
$ws.Range("I27").Value = 1759.3334
$ws.Range("J27").Value = 1802.8
$ws.Range("K27").Value = 1759.3334
$ws.Range("L27").Value = 1802.8
$ws.Range("M27").Value = -1652.3334
$ws.Range("N27").Value = -2016.8

$ws.Range("H46").Value = 1000
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").Value = ""

$ws.Range("H55").Value = 193.06667
$ws.Range("I55").Value = 164.2
$ws.Range("J55").Value = 250.8
$ws.Range("K55").Value = 164.2
$ws.Range("L55").Value = 250.8
$ws.Range("M55").Value = 8.800000000000011
$ws.Range("N55").Value = -596.8

$ws.Range("H126").Value = 6142.2666
$ws.Range("I126").Value = 5784.1113
$ws.Range("J126").Value = 6679.5
$ws.Range("K126").Value = 17352.3339
$ws.Range("L126").Value = 20038.5
$ws.Range("M126").Value = -14882.3339
$ws.Range("N126").Value = -24978.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2757.923
$ws.Range("I122").Value = 2441.182
$ws.Range("K122").Value = 7323.545999999999
$ws.Range("M122").Value = -4873.545999999999
